$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" (sheet1) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.9430604982206405
$summary.Range("C2").Value = 0.4285714285714285
$summary.Range("D2").Value = 0.4285714285714285
$summary.Range("E2").Value = 0.4285714285714285
$summary.Range("F2").Value = 0.4285714285714285
$summary.Range("G2").Value = 0.4285714285714285
$summary.Range("H2").Value = 0.6993044408774746
$summary.Range("I2").Value = 12
$summary.Range("J2").Value = 16
$summary.Range("K2").Value = 518
$summary.Range("L2").Value = 16

# --- Sheet "Classification Report" (sheet2) ---
$report = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$report.Range("B2").Value = 0.9700374531835206
$report.Range("C2").Value = 0.9700374531835206
$report.Range("D2").Value = 0.9700374531835206

# Row 3 - class "1"
$report.Range("B3").Value = 0.4285714285714285
$report.Range("C3").Value = 0.4285714285714285
$report.Range("D3").Value = 0.4285714285714285

# Row 4 - accuracy
$report.Range("B4").Value = 0.9430604982206405
$report.Range("C4").Value = 0.9430604982206405
$report.Range("D4").Value = 0.9430604982206405
$report.Range("E4").Value = 0.9430604982206405

# Row 5 - macro avg
$report.Range("B5").Value = 0.6993044408774746
$report.Range("C5").Value = 0.6993044408774746
$report.Range("D5").Value = 0.6993044408774746

# Row 6 - weighted avg
$report.Range("B6").Value = 0.9430604982206405
$report.Range("C6").Value = 0.9430604982206405
$report.Range("D6").Value = 0.9430604982206405

# --- Sheet "Confusion Matrix" (sheet3) ---
$confusion = $wb.Worksheets.Item("Confusion Matrix")
$confusion.Range("B2").Value = 518
$confusion.Range("C2").Value = 16
$confusion.Range("B3").Value = 16
$confusion.Range("C3").Value = 12
